$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-08 (row 21)
$ws.Range("B21").Value = 6216
$ws.Range("D21").Value = 5570927
$ws.Range("E21").Value = 896.2237773487774
$ws.Range("F21").Value = 7.897934386391259
$ws.Range("H21").Value = 27.16189466629779
